$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each updated cell as literal text (matches source data which stores
# numeric-looking price/volume figures as inline strings, not numbers).
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '57.013.22'
Set-TextValue 'D3' '2.396.49'
Set-TextValue 'E3' '  +1.82%  '
Set-TextValue 'E4' '  +0.12%  '
Set-TextValue 'D5' '505.13'
Set-TextValue 'E5' '  -0.93%  '
Set-TextValue 'D6' '132.46'
Set-TextValue 'E6' '  +4.33%  '
Set-TextValue 'D7' '0.999'
Set-TextValue 'E7' '  +0.05%  '
Set-TextValue 'E8' '  +0.38%  '
Set-TextValue 'D9' '2.410.14'
Set-TextValue 'E9' '  +1.79%  '
Set-TextValue 'D10' '0.0968'
Set-TextValue 'E10' '  +1.51%  '
Set-TextValue 'E11' '  -1.05%  '
Set-TextValue 'E12' '  +2.52%  '
Set-TextValue 'D13' '4.57'
Set-TextValue 'E13' '  -4.69%  '
Set-TextValue 'D14' '2.827.72'
Set-TextValue 'E14' '  +1.93%  '
Set-TextValue 'D15' '56.946.54'
Set-TextValue 'E15' '  +1.40%  '
Set-TextValue 'D16' '21.87'
Set-TextValue 'E16' '  +2.78%  '
Set-TextValue 'E17' '  +2.69%  '
Set-TextValue 'D18' '2.425.92'
Set-TextValue 'E18' '  +4.10%  '
Set-TextValue 'D19' '10.22'
Set-TextValue 'E19' '  +0.24%  '
Set-TextValue 'B20' 'BitcoinCash'
Set-TextValue 'C20' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D20' '309.85'
Set-TextValue 'E20' '  +0.16%  '
Set-TextValue 'B21' 'Polkadot'
Set-TextValue 'C21' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D21' '4.03'
Set-TextValue 'E21' '  -0.33%  '
Set-TextValue 'E22' '  +3.16%  '
Set-TextValue 'D23' '5.87'
Set-TextValue 'E23' '  +0.31%  '
Set-TextValue 'D24' '0.998'
Set-TextValue 'E24' '  +0.07%  '
Set-TextValue 'D25' '65.04'
Set-TextValue 'E25' '  +0.13%  '
Set-TextValue 'D26' '0.997'
Set-TextValue 'E26' '  -0.19%  '
Set-TextValue 'E27' '  +0.84%  '
Set-TextValue 'E28' '  -2.30%  '
Set-TextValue 'D29' '7.45'
Set-TextValue 'E29' '  +3.46%  '
Set-TextValue 'D30' '173.86'
Set-TextValue 'E30' '  -0.50%  '
Set-TextValue 'E31' '  +1.82%  '
Set-TextValue 'E32' '  +0.24%  '
Set-TextValue 'D33' '5.94'
Set-TextValue 'E33' '  -2.61%  '
Set-TextValue 'D34' '1.11'
Set-TextValue 'E34' '  +0.13%  '
Set-TextValue 'E35' '  +0.09%  '
Set-TextValue 'E36' '  +0.00%  '
Set-TextValue 'D37' '17.94'
Set-TextValue 'E37' '  +1.88%  '
Set-TextValue 'E38' '  +1.20%  '
Set-TextValue 'E39' '  +4.13%  '
Set-TextValue 'D40' '36.75'
Set-TextValue 'E40' '  +3.46%  '
Set-TextValue 'D41' '0.815'
Set-TextValue 'E41' '  +2.77%  '
Set-TextValue 'E42' '  +2.15%  '
Set-TextValue 'D43' '132.99'
Set-TextValue 'E43' '  +8.21%  '
Set-TextValue 'D44' '4.98'
Set-TextValue 'E45' '  +0.69%  '
Set-TextValue 'D46' '252.60'
Set-TextValue 'E46' '  +0.48%  '
Set-TextValue 'D47' '0.566'
Set-TextValue 'E47' '  +0.36%  '
Set-TextValue 'E48' '  +1.13%  '
Set-TextValue 'D49' '0.0489'
Set-TextValue 'E49' '  +1.17%  '
Set-TextValue 'D50' '0.0210'
Set-TextValue 'E50' '  +2.30%  '
Set-TextValue 'E51' '  +10.85%  '
